# Generate Report for Handoff
#
# The localization status report is regenerated: the workbook-wide status
# moves from "Handed back: in sync with en-US" to "Ready for handoff", and
# the handoff/handback timestamps recorded alongside it are refreshed. The
# "Status" column (narrower now that the new status text is shorter) is
# re-sized to fit the new content on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$newStatus = "Ready for handoff"
$wsOverview.Range("E2").Value = $newStatus   # zh-cn status column
$wsOverview.Range("F2").Value = $newStatus   # de-de status column
$wsZhCn.Range("C2").Value     = $newStatus
$wsDeDe.Range("C2").Value     = $newStatus

# --- Refresh timestamps recorded at report-generation time ---
# Overview!G2 and de-de!H2 share the "Latest Handback DateTime" value.
$wsOverview.Range("G2").Value = "2016-08-19 04:53:54"
$wsDeDe.Range("H2").Value     = "2016-08-19 04:53:54"
# zh-cn!H2 is the "Latest Handoff Datetime" value.
$wsZhCn.Range("H2").Value     = "2016-08-19 04:53:49"

# --- Re-fit the "Status" column now that its text is shorter ---
# Target display width (per the generated report) is ~17.216 characters.
# Excel's ColumnWidth setter snaps to its internal character-width grid, so
# the assigned value is pre-compensated to land on the closest achievable
# grid point to that target.
$statusColumnWidth = 16.333333333333336
$wsOverview.Range("E1").ColumnWidth = $statusColumnWidth
$wsOverview.Range("F1").ColumnWidth = $statusColumnWidth
$wsZhCn.Range("C1").ColumnWidth     = $statusColumnWidth
$wsDeDe.Range("C1").ColumnWidth     = $statusColumnWidth
